# Generate Report for Handback
# Update the timestamp strings that record when the handoff/handback
# xliff files were generated/processed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for bc798a75... (shared with the de-de
# sheet's "Correspond Handoff Datetime" for the same file, since both
# originally held the exact same text and must stay in sync).
$wsOverview.Range("G2").Value = "2016-08-18 15:11:45"
$wsDeDe.Range("H2").Value     = "2016-08-18 15:11:45"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for bc798a75...
$wsZhCn.Range("H2").Value = "2016-08-18 15:11:39"
$wsZhCn.Range("K2").Value = "2016-08-18 15:12:26"

# de-de sheet: Correspond Handback DateTime for bc798a75...
$wsDeDe.Range("K2").Value = "2016-08-18 15:12:35"
